$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix two mis-typed separators (comma -> period) in proveedor names ---
$ws.Range("E23").Value2 = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F23").Value2 = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E58").Value2 = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F58").Value2 = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E95").Value2 = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F95").Value2 = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E57").Value2 = "FERNANDEZ MARIO H. GALLICET OSCAR M"

# --- Reformat "Importe" column (H) amounts: comma decimal / period thousands -> plain period decimal ---
$importes = @{
    2 = "359.00"
    3 = "1532.00"
    4 = "715.00"
    5 = "35000.00"
    6 = "1019.56"
    7 = "3485.51"
    8 = "14178.75"
    9 = "107511.21"
    10 = "8118.50"
    11 = "4884.77"
    12 = "1807.15"
    13 = "1841.25"
    14 = "6772.40"
    15 = "4514.13"
    16 = "88.00"
    17 = "206.00"
    18 = "7039.50"
    19 = "1365.00"
    20 = "376.00"
    21 = "189.59"
    22 = "33.50"
    23 = "1883.30"
    24 = "5.70"
    25 = "225.14"
    26 = "57.37"
    27 = "19015.76"
    28 = "128.56"
    29 = "12.00"
    30 = "15.00"
    31 = "2296.92"
    32 = "1486.44"
    33 = "4712.00"
    34 = "2.10"
    35 = "44.32"
    36 = "9842.23"
    37 = "376.44"
    38 = "56.00"
    39 = "2961.40"
    40 = "21.20"
    41 = "592.98"
    42 = "4275.31"
    43 = "640.00"
    44 = "3815.00"
    45 = "1650.00"
    46 = "3374.00"
    47 = "236.00"
    48 = "5765.00"
    49 = "116.40"
    50 = "1234.00"
    51 = "835.00"
    52 = "9485.80"
    53 = "20332.00"
    54 = "33496.48"
    55 = "10225.00"
    56 = "1147.10"
    57 = "622.00"
    58 = "5442.44"
    59 = "740.00"
    60 = "748.00"
    61 = "1832.96"
    62 = "142.50"
    63 = "753.00"
    64 = "217.20"
    65 = "1066.00"
    66 = "4334.30"
    67 = "10508.00"
    68 = "1746.00"
    69 = "4438.55"
    70 = "104.55"
    71 = "74.00"
    72 = "7.00"
    73 = "1290.00"
    74 = "61.00"
    75 = "2541.50"
    76 = "750.00"
    77 = "824.24"
    78 = "291.00"
    79 = "996.60"
    80 = "255.96"
    81 = "46.60"
    82 = "299.85"
    83 = "200.00"
    84 = "500.00"
    85 = "1626.44"
    86 = "285.00"
    87 = "2541.50"
    88 = "1256.00"
    89 = "5085.21"
    90 = "600.00"
    91 = "700.00"
    92 = "120.00"
    93 = "12770.00"
    94 = "162.50"
    95 = "80.00"
    96 = "28.00"
    97 = "480.00"
    98 = "250.00"
    99 = "118.18"
    100 = "351.00"
    101 = "567.83"
    102 = "3304.56"
    103 = "7.50"
    104 = "18.90"
    105 = "125.26"
    106 = "78.62"
    107 = "300.00"
    108 = "2700.00"
    109 = "237.40"
    110 = "28.00"
    111 = "1050.00"
    112 = "21003.05"
    113 = "1217.44"
    114 = "1781.09"
    115 = "3515.66"
    116 = "6893.34"
    117 = "3735.00"
    118 = "631231.18"
    119 = "81000.00"
}

foreach ($row in $importes.Keys) {
    $cell = $ws.Cells.Item($row, 8)
    $cell.Value2 = "'" + $importes[$row]
}

# Clear the quote-prefix formatting picked up from the text-forcing apostrophe above
# so the cell style index matches the original (no explicit style).
foreach ($row in $importes.Keys) {
    $ws.Cells.Item($row, 8).Style = "Normal"
}